$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Donnée A1, voici des modifications. Il est 11h."

$ws.Range("G13").Select()
